$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "8 - Close Chat / Close greeting"
$ws.Range("A8").HorizontalAlignment = -4131

$ws.Range("B10").Select()
